$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet/tab to reflect the new "through" date
$ws.Name = "Through 2022-12-09"

# Update the row label for December to reflect the new "through" date
$ws.Range("A13").Value = "December (through 12-09)"

# Update the December row (row 13) values for columns C..I (B13 unchanged)
$ws.Range("C13").Value = 26
$ws.Range("D13").Value = 36
$ws.Range("E13").Value = 21
$ws.Range("F13").Value = 13
$ws.Range("G13").Value = 42
$ws.Range("H13").Value = 71
$ws.Range("I13").Value = 37

# Update the Total row (row 14) values for columns C..I (B14 unchanged)
$ws.Range("C14").Value = 589
$ws.Range("D14").Value = 857
$ws.Range("E14").Value = 703
$ws.Range("F14").Value = 547
$ws.Range("G14").Value = 1306
$ws.Range("H14").Value = 1714
$ws.Range("I14").Value = 1553

$wb.Save()
